# Scheduled-runner market-data refresh for Louisoix_Profits workbook.
# Updates cached numeric columns H:N (currentAveragePrice / currentAveragePriceNQ /
# currentAveragePriceHQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# on the affected rows of each job sheet (ALC, ARM, CRP, CUL, GSM, LTW, WVR) to the
# latest pulled market-board values. Columns A:G (recipe/leve metadata) are untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 33: H33, I33, K33, M33
$ws.Range("H33").Value = 280.47058
$ws.Range("I33").Value = 126.35714
$ws.Range("K33").Value = 126.35714
$ws.Range("M33").Value = 102.64286

# Row 76: H76, I76, J76, K76, L76, M76, N76
$ws.Range("H76").Value = 4158.3335
$ws.Range("I76").Value = 4158.3335
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 4158.3335
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -3843.3335
$ws.Range("N76").ClearContents()

# Row 79: H79, I79, J79, K79, L79, M79, N79
$ws.Range("H79").Value = 4158.3335
$ws.Range("I79").Value = 4158.3335
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 4158.3335
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -3066.3335
$ws.Range("N79").ClearContents()

# Row 98: H98, I98, J98, K98, L98, M98, N98
$ws.Range("H98").Value = 1617.1212
$ws.Range("I98").Value = 928.93335
$ws.Range("J98").Value = 8499
$ws.Range("K98").Value = 928.93335
$ws.Range("L98").Value = 8499
$ws.Range("M98").Value = 569.06665
$ws.Range("N98").Value = -11495

# Row 116: H116, I116, J116, K116, L116, M116, N116
$ws.Range("H116").Value = 117331.78
$ws.Range("I116").Value = 4958.8
$ws.Range("J116").Value = 257798
$ws.Range("K116").Value = 4958.8
$ws.Range("L116").Value = 257798
$ws.Range("M116").Value = -1516.8
$ws.Range("N116").Value = -264682

# Row 122: H122, I122, J122, K122, L122, M122, N122
$ws.Range("H122").Value = 1617.1212
$ws.Range("I122").Value = 928.93335
$ws.Range("J122").Value = 8499
$ws.Range("K122").Value = 2786.80005
$ws.Range("L122").Value = 25497
$ws.Range("M122").Value = -336.8000499999998
$ws.Range("N122").Value = -30397

# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 10394.4375
$ws.Range("I132").Value = 8665.786
$ws.Range("J132").Value = 22495
$ws.Range("K132").Value = 25997.358
$ws.Range("L132").Value = 67485
$ws.Range("M132").Value = -23467.358
$ws.Range("N132").Value = -72545

# Row 137: H137, I137, J137, K137, L137, M137, N137
$ws.Range("H137").Value = 1730.6765
$ws.Range("I137").Value = 1320.5385
$ws.Range("J137").Value = 1984.5714
$ws.Range("K137").Value = 3961.6155
$ws.Range("L137").Value = 5953.7142
$ws.Range("M137").Value = -1411.6155
$ws.Range("N137").Value = -11053.7142

# Row 138: H138, I138, J138, K138, L138, M138, N138
$ws.Range("H138").Value = 7163.52
$ws.Range("I138").Value = 9540
$ws.Range("J138").Value = 6710.857
$ws.Range("K138").Value = 28620
$ws.Range("L138").Value = 20132.571
$ws.Range("M138").Value = -23480
$ws.Range("N138").Value = -30412.571

# ---------------------------------------------------------------------------
# Sheet: ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 32: H32, I32, J32, K32, L32, M32, N32
$ws.Range("H32").Value = 13913.197
$ws.Range("I32").Value = 14637.567
$ws.Range("J32").Value = 9446.25
$ws.Range("K32").Value = 14637.567
$ws.Range("L32").Value = 9446.25
$ws.Range("M32").Value = -14350.567
$ws.Range("N32").Value = -10020.25

# Row 45: H45, I45, J45, K45, L45, M45, N45
$ws.Range("H45").Value = 3676.875
$ws.Range("I45").Value = 3344.5833
$ws.Range("J45").Value = 4673.75
$ws.Range("K45").Value = 3344.5833
$ws.Range("L45").Value = 4673.75
$ws.Range("M45").Value = -2967.5833
$ws.Range("N45").Value = -5427.75

# Row 74: H74, I74, J74, K74, L74, M74, N74
$ws.Range("H74").Value = 2454.825
$ws.Range("I74").Value = 2233.5833
$ws.Range("J74").Value = 4446
$ws.Range("K74").Value = 2233.5833
$ws.Range("L74").Value = 4446
$ws.Range("M74").Value = -1359.5833
$ws.Range("N74").Value = -6194

# Row 77: H77, I77, J77, K77, L77, M77, N77
$ws.Range("H77").Value = 2454.825
$ws.Range("I77").Value = 2233.5833
$ws.Range("J77").Value = 4446
$ws.Range("K77").Value = 11167.9165
$ws.Range("L77").Value = 22230
$ws.Range("M77").Value = -6799.916499999999
$ws.Range("N77").Value = -30966

# Row 126: H126, I126, K126, M126
$ws.Range("H126").Value = 15000
$ws.Range("I126").Value = 15000
$ws.Range("K126").Value = 45000
$ws.Range("M126").Value = -42530

# Row 132: H132, I132, K132, M132
$ws.Range("H132").Value = 48947.773
$ws.Range("I132").Value = 53342.7
$ws.Range("K132").Value = 160028.1
$ws.Range("M132").Value = -157498.1

# ---------------------------------------------------------------------------
# Sheet: CRP
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 31: H31, I31, J31, K31, L31, M31, N31
$ws.Range("H31").Value = 3572.8
$ws.Range("I31").Value = 4707
$ws.Range("J31").Value = 2962.077
$ws.Range("K31").Value = 4707
$ws.Range("L31").Value = 2962.077
$ws.Range("M31").Value = -4412
$ws.Range("N31").Value = -3552.077

# Row 34: H34, I34, J34, K34, L34, M34, N34
$ws.Range("H34").Value = 3572.8
$ws.Range("I34").Value = 4707
$ws.Range("J34").Value = 2962.077
$ws.Range("K34").Value = 4707
$ws.Range("L34").Value = 2962.077
$ws.Range("M34").Value = -4505
$ws.Range("N34").Value = -3366.077

# Row 58: H58, I58, J58, K58, L58, M58, N58
$ws.Range("H58").Value = 204300.2
$ws.Range("I58").Value = 335337
$ws.Range("J58").Value = 7745
$ws.Range("K58").Value = 335337
$ws.Range("L58").Value = 7745
$ws.Range("M58").Value = -335134
$ws.Range("N58").Value = -8151

# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 12000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 12000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 36000
$ws.Range("N132").Value = -41060
$ws.Range("M132").ClearContents()

# Row 136: H136, I136, J136, K136, L136, M136, N136
$ws.Range("H136").Value = 204300.2
$ws.Range("I136").Value = 335337
$ws.Range("J136").Value = 7745
$ws.Range("K136").Value = 1006011
$ws.Range("L136").Value = 23235
$ws.Range("M136").Value = -1003461
$ws.Range("N136").Value = -28335

# ---------------------------------------------------------------------------
# Sheet: CUL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 4: H4, I4, J4, K4, L4, M4, N4
$ws.Range("H4").Value = 559139
$ws.Range("I4").Value = 594053.9399999999
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 1782161.82
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -1782049.82
$ws.Range("N4").Value = -1724

# Row 5: H5, I5, J5, K5, L5, M5, N5
$ws.Range("H5").Value = 1535.4
$ws.Range("I5").Value = 800
$ws.Range("J5").Value = 1719.25
$ws.Range("K5").Value = 2400
$ws.Range("L5").Value = 5157.75
$ws.Range("M5").Value = -2288
$ws.Range("N5").Value = -5381.75

# Row 68: H68, I68, J68, K68, L68, M68, N68
$ws.Range("H68").Value = 2739.6
$ws.Range("I68").Value = 2349.5
$ws.Range("J68").Value = 3324.75
$ws.Range("K68").Value = 7048.5
$ws.Range("L68").Value = 9974.25
$ws.Range("M68").Value = -6237.5
$ws.Range("N68").Value = -11596.25

# Row 71: H71, I71, J71, K71, L71, M71, N71
$ws.Range("H71").Value = 2739.6
$ws.Range("I71").Value = 2349.5
$ws.Range("J71").Value = 3324.75
$ws.Range("K71").Value = 21145.5
$ws.Range("L71").Value = 29922.75
$ws.Range("M71").Value = -17089.5
$ws.Range("N71").Value = -38034.75

# Row 135: H135, I135, J135, K135, L135, M135, N135
$ws.Range("H135").Value = 1535.4
$ws.Range("I135").Value = 800
$ws.Range("J135").Value = 1719.25
$ws.Range("K135").Value = 7200
$ws.Range("L135").Value = 15473.25
$ws.Range("M135").Value = -4665
$ws.Range("N135").Value = -20543.25

# ---------------------------------------------------------------------------
# Sheet: GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Row 41: H41, I41, K41, M41
$ws.Range("H41").Value = 18492.5
$ws.Range("I41").Value = 16990
$ws.Range("K41").Value = 16990
$ws.Range("M41").Value = -16635

# Row 126: H126, I126, J126, K126, L126, M126, N126
$ws.Range("H126").Value = 6121.05
$ws.Range("I126").Value = 5133.077
$ws.Range("J126").Value = 7955.857
$ws.Range("K126").Value = 15399.231
$ws.Range("L126").Value = 23867.571
$ws.Range("M126").Value = -12929.231
$ws.Range("N126").Value = -28807.571

# Row 132: H132, I132, K132, M132
$ws.Range("H132").Value = 107599.7
$ws.Range("I132").Value = 107599.7
$ws.Range("K132").Value = 322799.1
$ws.Range("M132").Value = -320269.1

# ---------------------------------------------------------------------------
# Sheet: LTW
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 68: H68, I68, J68, K68, L68, M68, N68
$ws.Range("H68").Value = 7093.2
$ws.Range("I68").Value = 5000
$ws.Range("J68").Value = 7616.5
$ws.Range("K68").Value = 5000
$ws.Range("L68").Value = 7616.5
$ws.Range("M68").Value = -4251
$ws.Range("N68").Value = -9114.5

# Row 71: H71, I71, J71, K71, L71, M71, N71
$ws.Range("H71").Value = 7093.2
$ws.Range("I71").Value = 5000
$ws.Range("J71").Value = 7616.5
$ws.Range("K71").Value = 25000
$ws.Range("L71").Value = 38082.5
$ws.Range("M71").Value = -21256
$ws.Range("N71").Value = -45570.5

# Row 122: H122, I122, J122, K122, L122, M122, N122
$ws.Range("H122").Value = 4649
$ws.Range("I122").Value = 4470.4287
$ws.Range("J122").Value = 4732.3335
$ws.Range("K122").Value = 13411.2861
$ws.Range("L122").Value = 14197.0005
$ws.Range("M122").Value = -10961.2861
$ws.Range("N122").Value = -19097.0005

# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 40391.375
$ws.Range("I132").Value = 49750.48
$ws.Range("J132").Value = 6966
$ws.Range("K132").Value = 149251.44
$ws.Range("L132").Value = 20898
$ws.Range("M132").Value = -146721.44
$ws.Range("N132").Value = -25958

# ---------------------------------------------------------------------------
# Sheet: WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 81: H81, I81, K81, M81
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()

# Row 84: H84, I84, K84, M84
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()

# Row 132: H132, J132, L132, N132
$ws.Range("H132").Value = 36979.535
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
